# "Generate Report for Handback" - mark files as handed back (in sync with en-US)
# and record the handback xliff file + timestamp for each locale, plus point the
# "Latest Target File" column at the source file via a hyperlink.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$handedBackStatus = "Handed back: in sync with en-US"

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a95437b7ed69af95b67a93e7bad94afbf1bd960/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a95437b7ed69af95b67a93e7bad94afbf1bd960/e2e/b.md"

# ---------------------------------------------------------------------------
# Overview sheet: the per-locale status columns (E/F) move from "Ready for
# handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $handedBackStatus
$wsOverview.Range("F2").Value = $handedBackStatus
$wsOverview.Range("E3").Value = $handedBackStatus
$wsOverview.Range("F3").Value = $handedBackStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet: Status column, Latest Target File (+hyperlink), Latest
# Handback File and Latest Handback DateTime.
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $handedBackStatus
$wsZhCn.Range("C3").Value = $handedBackStatus

$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-02 14:46:10"
$wsZhCn.Range("K3").Value = "2016-09-02 14:46:10"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1

$wsZhCn.Hyperlinks.Delete()
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $aMdUrl, "", "", "a.md")
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $aMdUrl, "", "", "a.md")
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $bMdUrl, "", "", "b.md")
$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $aMdUrl, "", "", "a.md")

# ---------------------------------------------------------------------------
# de-de sheet: same shape of update, different handback file name/timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("C3").Value = $handedBackStatus

$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-02 14:46:24"
$wsDeDe.Range("K3").Value = "2016-09-02 14:46:24"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1

$wsDeDe.Hyperlinks.Delete()
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $aMdUrl, "", "", "a.md")
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $aMdUrl, "", "", "a.md")
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $bMdUrl, "", "", "b.md")
$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $aMdUrl, "", "", "a.md")
